$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144:166 down to 145:167.
# The new row formatting (e.g. the date style on column D) is inherited from the
# row above, which already matches the style used throughout this block.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new weekly record.
$ws.Cells.Item(144, 1).Value  = 4
$ws.Cells.Item(144, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(144, 3).Value  = "Los Lagos"
$ws.Cells.Item(144, 4).Value  = 44504
$ws.Cells.Item(144, 5).Value  = 10
$ws.Cells.Item(144, 6).Value  = "Fruta"
$ws.Cells.Item(144, 7).Value  = 100102
$ws.Cells.Item(144, 8).Value  = "Cítricos"
$ws.Cells.Item(144, 9).Value  = 100102006
$ws.Cells.Item(144, 10).Value = "Pomelo"
$ws.Cells.Item(144, 11).Value = "Start Ruby"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 120
$ws.Cells.Item(144, 14).Value = 11000
$ws.Cells.Item(144, 15).Value = 12000
$ws.Cells.Item(144, 16).Value = 11500
$ws.Cells.Item(144, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(144, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(144, 19).Value = 821
$ws.Cells.Item(144, 20).Value = 14
